$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locates the (1-based) index of the first paragraph whose visible text
# equals $oldText exactly (ignoring the trailing paragraph-mark
# character that Range.Text always carries).
function Find-ParagraphIndex($oldText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd()
        if ($t -eq $oldText) {
            return $i
        }
    }
    return -1
}

# Replaces the visible text of paragraph $paraIndex with $newText while
# preserving the paragraph's pPr, any leading empty <w:r/> run, and the
# rPr (bold/italic/etc.) of the run that carries the text - exactly what
# a plain text substitution should look like in the underlying XML.
function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range.WordOpenXML

    $pMatch = [regex]::Match($full, '<w:p[ >].*</w:p>')
    $pXml = $pMatch.Value

    $pPrXml = ""
    $pPrMatch = [regex]::Match($pXml, '<w:pPr>.*?</w:pPr>')
    if ($pPrMatch.Success) {
        $pPrXml = $pPrMatch.Value
    }

    $hasEmptyRun = $pXml.Contains("<w:r></w:r>") -or $pXml.Contains("<w:r/>")

    $rPrXml = ""
    $rPrMatch = [regex]::Match($pXml, '<w:r><w:rPr>.*?</w:rPr><w:t')
    if ($rPrMatch.Success) {
        $rPrXml = [regex]::Match($rPrMatch.Value, '<w:rPr>.*?</w:rPr>').Value
    }

    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $runsXml = ""
    if ($hasEmptyRun) {
        $runsXml = $runsXml + "<w:r/>"
    }
    $runsXml = $runsXml + "<w:r>" + $rPrXml + "<w:t>" + $escaped + "</w:t></w:r>"

    $newParaXml = "<w:p>" + $pPrXml + $runsXml + "</w:p>"
    $wrapXml = $pkgHeader + $newParaXml + $pkgFooter

    $isLast = ($paraIndex -eq $d.Paragraphs.Count)
    $countBefore = $d.Paragraphs.Count

    $rng = $p.Range
    $rng.MoveEnd(1, -1)
    $rng.Text = ""
    $rng.InsertXML($wrapXml)

    if ($isLast -and ($d.Paragraphs.Count -gt $countBefore)) {
        # Inserting a full <w:p> fragment right before the body's closing
        # sectPr mints an extra trailing empty paragraph (the old
        # paragraph mark survives as its own, now-empty, paragraph).
        # Merge it back away so the paragraph count returns to normal.
        $count = $d.Paragraphs.Count
        $prevEnd = $d.Paragraphs($count - 1).Range.End
        $lastEnd = $d.Paragraphs($count).Range.End
        $d.Range($prevEnd - 1, $lastEnd).Delete()
    }
}

# Finds the (first, not-yet-edited) paragraph matching $oldText and
# rewrites its text to $newText, preserving run/paragraph formatting.
function Replace-ParagraphText($oldText, $newText) {
    $idx = Find-ParagraphIndex $oldText
    if ($idx -lt 0) {
        Write-Output ("WARNING: paragraph not found for: " + $oldText)
        return
    }
    Set-ParagraphText $idx $newText
}

# Title - appears twice: once as the Heading1 and once as a bold run
# near the end of the document.
Replace-ParagraphText "Play Excalibur Unleashed for Free - Review" "Play Excalibur Unleashed for Free - Review & Game Details"
Replace-ParagraphText "Play Excalibur Unleashed for Free - Review" "Play Excalibur Unleashed for Free - Review & Game Details"

# "What we like" bullet list
Replace-ParagraphText "Unusual cross structure" "Unusual cross structure for a unique gaming experience"
Replace-ParagraphText "Impressive graphics and theme" "Symbols and theme related to the Arthurian legend"
Replace-ParagraphText "Simple gameplay for beginners" "Free Spin feature adds excitement to gameplay"
Replace-ParagraphText "Potential for high winnings" "Impressive graphics that immerse players in the medieval setting"

# "What we don't like" bullet list
Replace-ParagraphText "No bonus games or special features" "Lack of bonus games and special features"
Replace-ParagraphText "Not suitable for high rollers" "High volatility may not be suitable for high rollers"

# Meta description (italic run at the very end of the document)
Replace-ParagraphText "Read our review of Excalibur Unleashed slot and play it for free. Discover its unusual structure, special symbols, potential win, and overall gameplay experience." "Discover the magic of Excalibur Unleashed slot with its unique gameplay and impressive graphics. Play for free now!"
